$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Row 133 already had Waist/Wgt logged - fill in the kcal totals that were missing
$ws.Range("D133").Value = 2653
$ws.Range("E133").Value = 2326

# Rows 134-144 - newly logged bio / nutrition entries (Waist, Wgt, kcal Total, kcal)
$rows = @(
    @{ Row=134; B=96.5; C=79.2;               D=3501; E=3501 },
    @{ Row=135; B=96.5; C=79.3;               D=2782; E=2606 },
    @{ Row=136; B=97;   C=79.400000000000006; D=2364; E=1971 },
    @{ Row=137; B=97;   C=79.599999999999994; D=3450; E=2875 },
    @{ Row=138; B=97;   C=80.400000000000006; D=2135; E=2135 },
    @{ Row=139; B=97;   C=80.5;               D=2606; E=1997 },
    @{ Row=140; B=97;   C=79.599999999999994; D=2390; E=2027 },
    @{ Row=141; B=97;   C=79.7;               D=2747; E=2747 },
    @{ Row=142; B=96.5; C=79.599999999999994; D=2789; E=2789 },
    @{ Row=143; B=96.5; C=79.5;               D=2130; E=2130 },
    @{ Row=144; B=96.5; C=78.900000000000006; D=2060; E=1799 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B   # Waist (cm)
    $ws.Cells.Item($r.Row, 3).Value = $r.C   # Wgt (kg)
    $ws.Cells.Item($r.Row, 4).Value = $r.D   # kcal Total
    $ws.Cells.Item($r.Row, 5).Value = $r.E   # kcal
    $ws.Cells.Item($r.Row, 6).Value = 1      # Creatine
}

# Row 145 - only Waist/Wgt recorded so far, no kcal totals yet
$ws.Range("B145").Value = 96.5
$ws.Range("C145").Value = 78.8
$ws.Range("F145").Value = 1

# Rows 146-148 were placeholder future dates with no log entries yet - clear them back out
$ws.Range("A146").ClearContents()
$ws.Range("A147").ClearContents()
$ws.Range("A148").ClearContents()

# F134:F145 move from the "blank trailing row" border style to the normal
# filled-row style (matching F133), since they now hold real data
$ws.Range("F133").Copy()
$ws.Range("F134:F145").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Reflect where the user ended up scrolled to / selected
$excel.ActiveWindow.ScrollRow = 129
[void]$ws.Range("D142").Select()
